# Status.xlsx update:
#  - Append a 4th bullet point to the "Comments" cell (D13) on the
#    "Status(Summary)" sheet, describing the newly tested public access
#    methods.
#  - Row 13 grows (wrap-text) to fit the extra line, so bump its height.
#  - Move the sheet's active selection from D13 to B13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Status(Summary)")

$cell = $ws.Range("D13")
$existingText = $cell.Value2
$appendedText = "4) Public methods for getting a rendered survey and getting responses in flat tree done "
$cell.Value = $existingText + "`n" + $appendedText

# The extra wrapped line needs a taller row to stay fully visible.
$ws.Rows.Item(13).RowHeight = 143.35

# Author finished editing with B13 selected instead of D13.
$ws.Range("B13").Select() | Out-Null
